$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.406
$ws.Range("B3").Value = 6.701000000000001
$ws.Range("C5").Value = -13.066
$ws.Range("D5").Value = -8.021000000000001
$ws.Range("D9").Value = -7.634
$ws.Range("D11").Value = -8.25
$ws.Range("B14").Value = 6.449
$ws.Range("B21").Value = 6.054
$ws.Range("D21").Value = -7.775999999999999
$ws.Range("B23").Value = 6.444
$ws.Range("B25").Value = 6.103999999999999
